# Updates the crypto price/volume table to the latest scraped snapshot.
# Cell list is (A1 ref, new text) pairs taken verbatim from the source diff;
# numeric-looking "Price" values are written with a leading apostrophe so
# Excel keeps them as literal text (preserving things like trailing zeros
# in '1.010' or the two-dot thousands format in '26.938.12') instead of
# silently re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '26.938.12' }
    @{ Cell = 'D3'; Value = '1.845.79' }
    @{ Cell = 'E3'; Value = '  +0.94%  ' }
    @{ Cell = 'E4'; Value = '  +0.39%  ' }
    @{ Cell = 'D5'; Value = '309.66' }
    @{ Cell = 'E5'; Value = '  +0.36%  ' }
    @{ Cell = 'D6'; Value = '1.010' }
    @{ Cell = 'E6'; Value = '  +0.22%  ' }
    @{ Cell = 'D7'; Value = '0.4783' }
    @{ Cell = 'E7'; Value = '  +2.58%  ' }
    @{ Cell = 'D8'; Value = '0.3669' }
    @{ Cell = 'E8'; Value = '  +1.72%  ' }
    @{ Cell = 'D9'; Value = '0.07225' }
    @{ Cell = 'E9'; Value = '  +1.12%  ' }
    @{ Cell = 'D10'; Value = '0.9265' }
    @{ Cell = 'E10'; Value = '  +2.35%  ' }
    @{ Cell = 'D11'; Value = '19.68' }
    @{ Cell = 'E11'; Value = '  +1.22%  ' }
    @{ Cell = 'D12'; Value = '0.07708' }
    @{ Cell = 'E12'; Value = '  -0.82%  ' }
    @{ Cell = 'D13'; Value = '1.896.67' }
    @{ Cell = 'E13'; Value = '  +3.60%  ' }
    @{ Cell = 'E14'; Value = '  +0.87%  ' }
    @{ Cell = 'D15'; Value = '6.408' }
    @{ Cell = 'E15'; Value = '  +0.86%  ' }
    @{ Cell = 'D16'; Value = '88.81' }
    @{ Cell = 'E16'; Value = '  +1.06%  ' }
    @{ Cell = 'E17'; Value = '  +0.33%  ' }
    @{ Cell = 'E18'; Value = '  +0.84%  ' }
    @{ Cell = 'D19'; Value = '1.010' }
    @{ Cell = 'E19'; Value = '  +0.26%  ' }
    @{ Cell = 'D20'; Value = '26.972.37' }
    @{ Cell = 'E20'; Value = '  +0.87%  ' }
    @{ Cell = 'D21'; Value = '14.55' }
    @{ Cell = 'E21'; Value = '  +2.34%  ' }
    @{ Cell = 'E22'; Value = '  +0.88%  ' }
    @{ Cell = 'E23'; Value = '  +0.88%  ' }
    @{ Cell = 'D24'; Value = '1.928' }
    @{ Cell = 'E24'; Value = '  +0.31%  ' }
    @{ Cell = 'D25'; Value = '152.44' }
    @{ Cell = 'D26'; Value = '18.16' }
    @{ Cell = 'E26'; Value = '  +1.26%  ' }
    @{ Cell = 'D27'; Value = '1.994' }
    @{ Cell = 'E27'; Value = '  +0.85%  ' }
    @{ Cell = 'D28'; Value = '114.11' }
    @{ Cell = 'E28'; Value = '  +0.21%  ' }
    @{ Cell = 'D29'; Value = '4.943' }
    @{ Cell = 'E29'; Value = '  +2.28%  ' }
    @{ Cell = 'D30'; Value = '0.08885' }
    @{ Cell = 'E30'; Value = '  +0.86%  ' }
    @{ Cell = 'D31'; Value = '3.317' }
    @{ Cell = 'E31'; Value = '  +5.32%  ' }
    @{ Cell = 'D32'; Value = '1.172' }
    @{ Cell = 'E32'; Value = '  +1.85%  ' }
    @{ Cell = 'D33'; Value = '0.7440' }
    @{ Cell = 'E33'; Value = '  +1.38%  ' }
    @{ Cell = 'D34'; Value = '4.492' }
    @{ Cell = 'E34'; Value = '  +0.90%  ' }
    @{ Cell = 'D35'; Value = '2.719' }
    @{ Cell = 'E35'; Value = '  -0.71%  ' }
    @{ Cell = 'D36'; Value = '1.128' }
    @{ Cell = 'E36'; Value = '  +4.42%  ' }
    @{ Cell = 'D37'; Value = '0.01956' }
    @{ Cell = 'E37'; Value = '  +1.60%  ' }
    @{ Cell = 'D38'; Value = '0.05265' }
    @{ Cell = 'E38'; Value = '  +2.58%  ' }
    @{ Cell = 'D39'; Value = '2.983' }
    @{ Cell = 'E39'; Value = '  +1.76%  ' }
    @{ Cell = 'D40'; Value = '0.5194' }
    @{ Cell = 'E40'; Value = '  +2.29%  ' }
    @{ Cell = 'D41'; Value = '6.994' }
    @{ Cell = 'E41'; Value = '  +1.56%  ' }
    @{ Cell = 'D42'; Value = '0.1510' }
    @{ Cell = 'E42'; Value = '  +0.76%  ' }
    @{ Cell = 'D43'; Value = '8.200' }
    @{ Cell = 'E43'; Value = '  +1.76%  ' }
    @{ Cell = 'D44'; Value = '10.53' }
    @{ Cell = 'E44'; Value = '  +4.89%  ' }
    @{ Cell = 'D45'; Value = '0.4727' }
    @{ Cell = 'E45'; Value = '  +1.07%  ' }
    @{ Cell = 'E46'; Value = '  +0.31%  ' }
    @{ Cell = 'D47'; Value = '101.45' }
    @{ Cell = 'E47'; Value = '  +3.07%  ' }
    @{ Cell = 'E48'; Value = '  +2.39%  ' }
    @{ Cell = 'D49'; Value = '65.53' }
    @{ Cell = 'E49'; Value = '  +2.20%  ' }
    @{ Cell = 'D50'; Value = '0.06029' }
    @{ Cell = 'E50'; Value = '  -0.36%  ' }
    @{ Cell = 'D51'; Value = '0.8858' }
    @{ Cell = 'E51'; Value = '  +3.74%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $text = $u.Value
    $looksNumeric = $text -match '^[+-]?(\d+\.?\d*|\.\d+)([eE][+-]?\d+)?$'
    if ($looksNumeric) {
        # Prefix with an apostrophe (Excel's own force-text marker) so the
        # value round-trips as text, then restore the default style so the
        # only lasting change is the cell's content.
        $range.Value = "'" + $text
        $range.Style = 'Normal'
    } else {
        $range.Value = $text
    }
}
